$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1157.3462
$ws.Range("I19").Value = 202.13333
$ws.Range("J19").Value = 2459.9092
$ws.Range("K19").Value = 202.13333
$ws.Range("L19").Value = 2459.9092
$ws.Range("M19").Value = -27.13333
$ws.Range("N19").Value = -2809.9092
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H40").Value = 2500.8333
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2601
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2601
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2951
$ws.Range("H86").Value = 214914.28
$ws.Range("I86").Value = 899.3333
$ws.Range("K86").Value = 899.3333
$ws.Range("M86").Value = 223.6667
$ws.Range("H89").Value = 214914.28
$ws.Range("I89").Value = 899.3333
$ws.Range("K89").Value = 4496.6665
$ws.Range("M89").Value = 1119.3335
$ws.Range("H94").Value = 10216.889
$ws.Range("I94").Value = 10216.889
$ws.Range("K94").Value = 10216.889
$ws.Range("M94").Value = -9765.888999999999
$ws.Range("H100").Value = 1082.4166
$ws.Range("I100").Value = 1122.4445
$ws.Range("J100").Value = 962.3333
$ws.Range("K100").Value = 1122.4445
$ws.Range("L100").Value = 962.3333
$ws.Range("M100").Value = -581.4445000000001
$ws.Range("N100").Value = -2044.3333
$ws.Range("H107").Value = 37038292
$ws.Range("I107").Value = 45455892
$ws.Range("K107").Value = 45455892
$ws.Range("M107").Value = -45453972
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H127").Value = 4597
$ws.Range("I127").Value = 4597
$ws.Range("K127").Value = 13791
$ws.Range("M127").Value = -8831
$ws.Range("H132").Value = 2230.037
$ws.Range("I132").Value = 1273.0416
$ws.Range("K132").Value = 3819.1248
$ws.Range("M132").Value = -1289.1248
$ws.Range("H138").Value = 3069.3572
$ws.Range("J138").Value = 3405.0312
$ws.Range("L138").Value = 10215.0936
$ws.Range("N138").Value = -20495.0936
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 829.1667
$ws.Range("I5").Value = 595
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 595
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -483
$ws.Range("N5").Value = -2224
$ws.Range("H61").Value = 7841.6
$ws.Range("I61").Value = 7841.6
$ws.Range("K61").Value = 7841.6
$ws.Range("M61").Value = -7629.6
$ws.Range("H74").Value = 3623
$ws.Range("J74").Value = 3997.3333
$ws.Range("L74").Value = 3997.3333
$ws.Range("N74").Value = -5745.3333
$ws.Range("H77").Value = 3623
$ws.Range("J77").Value = 3997.3333
$ws.Range("L77").Value = 19986.6665
$ws.Range("N77").Value = -28722.6665
$ws.Range("H132").Value = 5342.6665
$ws.Range("I132").Value = 4214.6
$ws.Range("K132").Value = 12643.8
$ws.Range("M132").Value = -10113.8
$ws.Range("H136").Value = 7841.6
$ws.Range("I136").Value = 7841.6
$ws.Range("K136").Value = 23524.8
$ws.Range("M136").Value = -20974.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 829.1667
$ws.Range("I4").Value = 595
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 595
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -480
$ws.Range("N4").Value = -2230
$ws.Range("H17").Value = 2900
$ws.Range("J17").Value = 2900
$ws.Range("L17").Value = 2900
$ws.Range("N17").Value = -3244
$ws.Range("H64").Value = 1249.6666
$ws.Range("I64").Value = 1125
$ws.Range("J64").Value = 1499
$ws.Range("K64").Value = 1125
$ws.Range("L64").Value = 1499
$ws.Range("M64").Value = -900
$ws.Range("N64").Value = -1949
$ws.Range("H67").Value = 1249.6666
$ws.Range("I67").Value = 1125
$ws.Range("J67").Value = 1499
$ws.Range("K67").Value = 1125
$ws.Range("L67").Value = 1499
$ws.Range("M67").Value = -345
$ws.Range("N67").Value = -3059
$ws.Range("H107").Value = 2587.25
$ws.Range("I107").Value = 2449.6667
$ws.Range("K107").Value = 2449.6667
$ws.Range("M107").Value = -529.6667000000002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 248.625
$ws.Range("I15").Value = 150
$ws.Range("J15").Value = 347.25
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = 347.25
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -687.25
$ws.Range("H16").Value = 83336590
$ws.Range("I16").Value = 83336590
$ws.Range("K16").Value = 83336590
$ws.Range("M16").Value = -83336303
$ws.Range("H22").Value = 67566.89999999999
$ws.Range("I22").Value = 86274.07000000001
$ws.Range("J22").Value = 23916.834
$ws.Range("K22").Value = 86274.07000000001
$ws.Range("L22").Value = 23916.834
$ws.Range("M22").Value = -85924.07000000001
$ws.Range("N22").Value = -24616.834
$ws.Range("H29").Value = 28333
$ws.Range("I29").Value = 28000
$ws.Range("J29").Value = 28499.5
$ws.Range("K29").Value = 28000
$ws.Range("L29").Value = 28499.5
$ws.Range("M29").Value = -27707
$ws.Range("N29").Value = -29085.5
$ws.Range("H113").Value = 83336590
$ws.Range("I113").Value = 83336590
$ws.Range("K113").Value = 83336590
$ws.Range("M113").Value = -83334420
$ws.Range("H122").Value = 1661.4
$ws.Range("I122").Value = 1214.25
$ws.Range("J122").Value = 3450
$ws.Range("K122").Value = 3642.75
$ws.Range("L122").Value = 10350
$ws.Range("M122").Value = -1192.75
$ws.Range("N122").Value = -15250
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2285.4285
$ws.Range("I98").Value = 2259.6
$ws.Range("J98").Value = 2350
$ws.Range("K98").Value = 6778.799999999999
$ws.Range("L98").Value = 7050
$ws.Range("M98").Value = -5280.799999999999
$ws.Range("N98").Value = -10046
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 27006250
$ws.Range("I11").Value = 38000000
$ws.Range("J11").Value = 8683333
$ws.Range("K11").Value = 38000000
$ws.Range("L11").Value = 8683333
$ws.Range("M11").Value = -37999861
$ws.Range("N11").Value = -8683611
$ws.Range("H107").Value = 1720.4286
$ws.Range("I107").Value = 338.5
$ws.Range("J107").Value = 3563
$ws.Range("K107").Value = 338.5
$ws.Range("L107").Value = 3563
$ws.Range("M107").Value = 1581.5
$ws.Range("N107").Value = -7403
$ws.Range("H132").Value = 4483.722
$ws.Range("I132").Value = 4582.8823
$ws.Range("K132").Value = 13748.6469
$ws.Range("M132").Value = -11218.6469
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2746.9092
$ws.Range("I46").Value = 1933.75
$ws.Range("K46").Value = 1933.75
$ws.Range("M46").Value = -1745.75
$ws.Range("H82").Value = 92183
$ws.Range("I82").Value = 1189.6666
$ws.Range("K82").Value = 1189.6666
$ws.Range("M82").Value = -828.6666
$ws.Range("H85").Value = 92183
$ws.Range("I85").Value = 1189.6666
$ws.Range("K85").Value = 1189.6666
$ws.Range("M85").Value = 58.33339999999998
$ws.Range("H93").Value = 835.6
$ws.Range("I93").Value = 835.6
$ws.Range("K93").Value = 835.6
$ws.Range("M93").Value = 412.4
$ws.Range("H100").Value = 7198.6
$ws.Range("I100").Value = 6998.5
$ws.Range("K100").Value = 6998.5
$ws.Range("M100").Value = -6457.5
$ws.Range("H122").Value = 5764.1875
$ws.Range("I122").Value = 4529
$ws.Range("J122").Value = 6999.375
$ws.Range("K122").Value = 13587
$ws.Range("L122").Value = 20998.125
$ws.Range("M122").Value = -11137
$ws.Range("N122").Value = -25898.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 20838.334
$ws.Range("I22").Value = 11250
$ws.Range("J22").Value = 40015
$ws.Range("K22").Value = 11250
$ws.Range("L22").Value = 40015
$ws.Range("M22").Value = -10957
$ws.Range("N22").Value = -40601
$ws.Range("H47").Value = 16799.2
$ws.Range("I47").Value = 9000
$ws.Range("J47").Value = 18749
$ws.Range("K47").Value = 9000
$ws.Range("L47").Value = 18749
$ws.Range("M47").Value = -8428
$ws.Range("N47").Value = -19893
$ws.Range("H107").Value = 299.58334
$ws.Range("I107").Value = 317.72726
$ws.Range("K107").Value = 953.18178
$ws.Range("M107").Value = 966.81822
$ws.Range("H113").Value = 1593.8889
$ws.Range("I113").Value = 407.66666
$ws.Range("J113").Value = 3966.3333
$ws.Range("K113").Value = 1222.99998
$ws.Range("L113").Value = 11898.9999
$ws.Range("M113").Value = 947.0000199999999
$ws.Range("N113").Value = -16238.9999
